$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

$ws.Cells.Item($row, 1).Value = 45744.16666666666
$ws.Cells.Item($row, 2).Value = 565.53
$ws.Cells.Item($row, 3).Value = 566.2675
$ws.Cells.Item($row, 4).Value = 555.0700000000001
$ws.Cells.Item($row, 5).Value = 555.66
$ws.Cells.Item($row, 6).Value = 71628953
$ws.Cells.Item($row, 7).Value = 516430404000
$ws.Cells.Item($row, 8).Value = 566.564
$ws.Cells.Item($row, 9).Value = 567.3779999999999
$ws.Cells.Item($row, 10).Value = 578.5740000000001
$ws.Cells.Item($row, 11).Value = 588.4064000000001
$ws.Cells.Item($row, 12).Value = ""
$ws.Cells.Item($row, 13).Value = -6.537610499655671
$ws.Cells.Item($row, 14).Value = -7.313184492112273
$ws.Cells.Item($row, 15).Value = 42.99157028033723
$ws.Cells.Item($row, 16).Value = 35.17092700954727
$ws.Cells.Item($row, 17).Value = 34.1669106881405
$ws.Cells.Item($row, 18).Value = 42.44607244607244
$ws.Cells.Item($row, 19).Value = ""
